$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full data set (header + all VSAC code system -> FHIR URI rows) reflecting
# the updated, alphabetically-ordered mapping table (11 rows added, "SOP" renamed to "SOPT").
$data = @(
    @("VSAC Code System Name", "FHIR URI"),
    @("ActCode", "http://terminology.hl7.org/CodeSystem/v3-ActCode"),
    @("ActMood", "http://terminology.hl7.org/CodeSystem/v3-ActMood"),
    @("ActPriority", "http://terminology.hl7.org/CodeSystem/v3-ActPriority"),
    @("ActReason", "http://terminology.hl7.org/CodeSystem/v3-ActReason"),
    @("ActRelationshipType", "http://terminology.hl7.org/CodeSystem/v3-ActRelationshipType"),
    @("ActStatus", "http://terminology.hl7.org/CodeSystem/v3-ActStatus"),
    @("AddressUse", "http://terminology.hl7.org/CodeSystem/v3-AddressUse"),
    @("AdministrativeGender", "http://terminology.hl7.org/CodeSystem/v3-AdministrativeGender"),
    @("AdministrativeSex", "http://terminology.hl7.org/CodeSystem/v2-0001"),
    @("AllergyIntoleranceClinicalStatusCodes", "http://terminology.hl7.org/CodeSystem/allergyintolerance-clinical"),
    @("AllergyIntoleranceVerificationStatus", "http://terminology.hl7.org/CodeSystem/allergyintolerance-verification"),
    @("CDCNHSN", "https://www.cdc.gov/nhsn/cdaportal/terminology/codesystem/cdcnhsn.html"),
    @("CDCREC", "urn:oid:2.16.840.1.113883.6.238"),
    @("CDT", "http://www.ada.org/cdt"),
    @("CPT", "http://www.ama-assn.org/go/cpt"),
    @("CVX", "http://hl7.org/fhir/sid/cvx"),
    @("ConditionClinicalStatusCodes", "http://terminology.hl7.org/CodeSystem/condition-clinical"),
    @("ConditionVerificationStatus", "http://terminology.hl7.org/CodeSystem/condition-ver-status"),
    @("Confidentiality", "http://terminology.hl7.org/CodeSystem/v3-Confidentiality"),
    @("DataAbsentReason", "http://terminology.hl7.org/CodeSystem/data-absent-reason"),
    @("DiagnosisRole", "http://terminology.hl7.org/CodeSystem/diagnosis-role"),
    @("DischargeDisposition", "http://terminology.hl7.org/CodeSystem/v2-0112"),
    @("EntityNamePartQualifier", "http://terminology.hl7.org/CodeSystem/v3-EntityNamePartQualifier"),
    @("EntityNameUse", "http://terminology.hl7.org/CodeSystem/v3-EntityNameUse"),
    @("EventTiming", "http://hl7.org/fhir/event-timing"),
    @("ExampleDiagnosisTypeCodes", "http://terminology.hl7.org/CodeSystem/ex-diagnosistype"),
    @("ExampleProcedureTypeCodes", "http://terminology.hl7.org/CodeSystem/ex-procedure-type"),
    @("HCPCS", "http://www.cms.gov/Medicare/Coding/HCPCSReleaseCodeSets"),
    @("HSLOC", "https://www.cdc.gov/nhsn/cdaportal/terminology/codesystem/hsloc.html"),
    @("HealthCareProviderTaxonomy", "http://nucc.org/provider-taxonomy"),
    @("ICD10CM", "http://hl7.org/fhir/sid/icd-10-cm"),
    @("ICD10PCS", "http://www.cms.gov/Medicare/Coding/ICD10"),
    @("ICD9CM", "http://hl7.org/fhir/sid/icd-9-cm"),
    @("LOINC", "http://loinc.org"),
    @("LanguageAbilityMode", "http://terminology.hl7.org/CodeSystem/v3-LanguageAbilityMode"),
    @("LanguageAbilityProficiency", "http://terminology.hl7.org/CodeSystem/v3-LanguageAbilityProficiency"),
    @("LivingArrangement", "http://terminology.hl7.org/CodeSystem/v3-LivingArrangement"),
    @("MED-RT", "http://www.nlm.nih.gov/research/umls/MED-RT"),
    @("MaritalStatus", "http://terminology.hl7.org/CodeSystem/v3-MaritalStatus"),
    @("MediaType", "http://terminology.hl7.org/CodeSystem/v3-mediaType"),
    @("MedicationRequestCategoryCodes", "http://terminology.hl7.org/CodeSystem/medicationrequest-category"),
    @("NCI", "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl"),
    @("NDFRT", "http://hl7.org/fhir/ndfrt"),
    @("NullFlavor", "http://terminology.hl7.org/CodeSystem/v3-NullFlavor"),
    @("ObservationCategoryCodes", "http://terminology.hl7.org/CodeSystem/observation-category"),
    @("ObservationInterpretation", "http://terminology.hl7.org/CodeSystem/v3-ObservationInterpretation"),
    @("ObservationValue", "http://terminology.hl7.org/CodeSystem/v3-ObservationValue"),
    @("ParticipationFunction", "http://terminology.hl7.org/CodeSystem/v3-ParticipationFunction"),
    @("ParticipationMode", "http://terminology.hl7.org/CodeSystem/v3-ParticipationMode"),
    @("ParticipationType", "http://terminology.hl7.org/CodeSystem/v3-ParticipationType"),
    @("PresentOnAdmission", "https://www.cms.gov/Medicare/Medicare-Fee-for-Service-Payment/HospitalAcqCond/Coding"),
    @("RXNORM", "http://www.nlm.nih.gov/research/umls/rxnorm"),
    @("ReligiousAffiliation", "http://terminology.hl7.org/CodeSystem/v3-ReligiousAffiliation"),
    @("RoleClass", "http://terminology.hl7.org/CodeSystem/v3-RoleClass"),
    @("RoleCode", "http://terminology.hl7.org/CodeSystem/v3-RoleCode"),
    @("RoleStatus", "http://terminology.hl7.org/CodeSystem/v3-RoleStatus"),
    @("SNOMEDCT", "http://snomed.info/sct"),
    @("SOPT", "https://nahdo.org/sopt"),
    @("TaskCode", "http://hl7.org/fhir/CodeSystem/task-code"),
    @("TimingEvent", "http://terminology.hl7.org/CodeSystem/v3-TimingEvent"),
    @("UCUM", "http://unitsofmeasure.org"),
    @("UMLS", "http://www.nlm.nih.gov/research/umls"),
    @("UNII", "http://fdasis.nlm.nih.gov"),
    @("USCoreConditionCategoryExtensionCodes", "http://hl7.org/fhir/us/core/CodeSystem/condition-category")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $ws.Cells.Item($i + 1, 1).Value = $row[0]
    $ws.Cells.Item($i + 1, 2).Value = $row[1]
}

$ws.Range("A1").Select() | Out-Null
